$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / header updates (Volume/Issue number and date range)
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# Row 14: Murder
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 8
$ws.Range("K14").Value = -25
$ws.Range("L14").Value = 50
$ws.Range("M14").Value = -40
$ws.Range("N14").Value = -57.142857142857

# Row 15: Rape
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("I15").Value = 17
$ws.Range("J15").Value = 23
$ws.Range("K15").Value = -26.086956521739
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -48.484848484848

# Row 16: Robbery
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = 8.823529411764
$ws.Range("I16").Value = 189
$ws.Range("J16").Value = 170
$ws.Range("K16").Value = 11.176470588235
$ws.Range("L16").Value = 34.042553191489
$ws.Range("M16").Value = 8
$ws.Range("N16").Value = -70.32967032967

# Row 17: Fel. Assault
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = -5
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = 13.043478260869
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 306
$ws.Range("K17").Value = 14.37908496732
$ws.Range("L17").Value = 19.047619047619
$ws.Range("M17").Value = 85.185185185185
$ws.Range("N17").Value = -5.149051490514

# Row 18: Burglary
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -13.636363636363
$ws.Range("I18").Value = 129
$ws.Range("J18").Value = 127
$ws.Range("K18").Value = 1.574803149606
$ws.Range("L18").Value = 27.722772277227
$ws.Range("M18").Value = -14
$ws.Range("N18").Value = -85.407239819004

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -23.529411764705
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = -29.850746268656
$ws.Range("I19").Value = 328
$ws.Range("J19").Value = 351
$ws.Range("K19").Value = -6.552706552706
$ws.Range("L19").Value = 58.454106280193
$ws.Range("M19").Value = 162.4
$ws.Range("N19").Value = 31.726907630522

# Row 20: G.L.A.
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 166.666666666667
$ws.Range("F20").Value = 37
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 68.181818181818
$ws.Range("I20").Value = 287
$ws.Range("J20").Value = 222
$ws.Range("K20").Value = 29.279279279279
$ws.Range("L20").Value = 23.706896551724
$ws.Range("M20").Value = 99.305555555555
$ws.Range("N20").Value = -59.97210599721

# Row 21: TOTAL
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = 1.88679245283
$ws.Range("F21").Value = 220
$ws.Range("G21").Value = 219
$ws.Range("H21").Value = 0.456621004566
$ws.Range("I21").Value = 1306
$ws.Range("J21").Value = 1207
$ws.Range("K21").Value = 8.202154101077
$ws.Range("L21").Value = 31.124497991967
$ws.Range("M21").Value = 61.234567901234
$ws.Range("N21").Value = -55.012056493282

# Row 22: Transit
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = -61.111111111111
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -58.823529411764
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = -62.5
$ws.Range("I23").Value = 44
$ws.Range("J23").Value = 49
$ws.Range("K23").Value = -10.204081632653
$ws.Range("L23").Value = -2.222222222222
$ws.Range("M23").Value = 76
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = -21.428571428571
$ws.Range("F24").Value = 127
$ws.Range("G24").Value = 121
$ws.Range("H24").Value = 4.95867768595
$ws.Range("I24").Value = 605
$ws.Range("J24").Value = 686
$ws.Range("K24").Value = -11.807580174927
$ws.Range("L24").Value = 33.849557522123
$ws.Range("M24").Value = 113.028169014085
$ws.Range("N24").Value = "***.*"

# Row 25: Misd. Assault
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 26.666666666666
$ws.Range("F25").Value = 86
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = 17.808219178082
$ws.Range("I25").Value = 416
$ws.Range("J25").Value = 392
$ws.Range("K25").Value = 6.122448979591
$ws.Range("L25").Value = 18.518518518518
$ws.Range("M25").Value = -5.022831050228
$ws.Range("N25").Value = "***.*"

# Row 26: UCR Rape*
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -66.666666666666
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = -70
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 41
$ws.Range("K26").Value = -36.585365853658
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 28.571428571428
$ws.Range("I27").Value = 39
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = 44.444444444444
$ws.Range("L27").Value = 18.181818181818
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28: Shooting Vic.
$ws.Range("C28").Value = "0"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = -25
$ws.Range("L28").Value = -25
$ws.Range("M28").Value = -30.76923076923
$ws.Range("N28").Value = -66.037735849056

# Row 29: Shooting Inc.
$ws.Range("C29").Value = "0"
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 15
$ws.Range("J29").Value = 23
$ws.Range("K29").Value = -34.782608695652
$ws.Range("L29").Value = -28.571428571428
$ws.Range("M29").Value = -34.782608695652
$ws.Range("N29").Value = -69.38775510204

# Row 30: Hate Crimes
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = "0"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("I30").Value = "0"
$ws.Range("J30").Value = 1
$ws.Range("K30").Value = -100
$ws.Range("L30").Value = "***.*"
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"

